$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells AD1:AF1 should look exactly like the existing header
# cells (bold, centered, bordered) -> copy the format from AC1 first,
# then fill in the text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-37: team win/loss/tie record (same values for every row).
for ($r = 2; $r -le 37; $r++) {
    $ws.Cells.Item($r, 30).Value = 95   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 67   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
